# Regenerate the "K" column (col G) save_data values for festa_matt.
# These values are computed/simulated upstream (s_vals) and then written
# directly into the sheet, replacing the previous Strike# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
